# Updated cryptos list on Mon Jun 10 03:37:00 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cellRef, $textValue)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $textValue
    $rng.Style = "Normal"
}

Set-TextCell "D2" "69.744.92"
Set-TextCell "E2" "  +0.79%  "

Set-TextCell "D3" "3.691.29"

Set-TextCell "D4" "1.00"
Set-TextCell "E4" "  +0.09%  "

Set-TextCell "D5" "672.04"
Set-TextCell "E5" "  +0.12%  "

Set-TextCell "D6" "160.61"
Set-TextCell "E6" "  +2.27%  "

Set-TextCell "E7" "  -0.01%  "

Set-TextCell "E8" "  +1.63%  "

Set-TextCell "D9" "0.146"
Set-TextCell "E9" "  +0.62%  "

Set-TextCell "E10" "  +1.92%  "

Set-TextCell "D11" "0.443"
Set-TextCell "E11" "  +1.81%  "

Set-TextCell "E12" "  +1.38%  "

Set-TextCell "D13" "33.17"
Set-TextCell "E13" "  +3.17%  "

Set-TextCell "D14" "3.670.73"
Set-TextCell "E14" "  +0.02%  "

Set-TextCell "D15" "69.720.10"

Set-TextCell "E16" "  +2.52%  "

Set-TextCell "D17" "16.16"
Set-TextCell "E17" "  +1.01%  "

Set-TextCell "D18" "6.49"
Set-TextCell "E18" "  +1.36%  "

Set-TextCell "D19" "471.19"
Set-TextCell "E19" "  +0.93%  "

Set-TextCell "D20" "9.78"
Set-TextCell "E20" "  -1.53%  "

Set-TextCell "E21" "  +0.19%  "

Set-TextCell "D22" "79.90"
Set-TextCell "E22" "  +0.40%  "

Set-TextCell "D23" "3.837.89"
Set-TextCell "E23" "  +0.57%  "

Set-TextCell "E24" "  +5.95%  "

Set-TextCell "E25" "  -0.01%  "

Set-TextCell "E26" "  +0.98%  "

Set-TextCell "E27" "  +0.87%  "

Set-TextCell "E28" "  +1.50%  "

Set-TextCell "E29" "  -0.55%  "

Set-TextCell "D30" "2.01"
Set-TextCell "E30" "  +2.08%  "

Set-TextCell "D31" "0.167"
Set-TextCell "E31" "  +4.91%  "

Set-TextCell "E32" "  +0.06%  "

Set-TextCell "D33" "26.84"
Set-TextCell "E33" "  +0.25%  "

Set-TextCell "E34" "  -1.30%  "

Set-TextCell "D35" "3.687.39"
Set-TextCell "E35" "  +0.63%  "

Set-TextCell "D36" "8.50"
Set-TextCell "E36" "  +4.86%  "

Set-TextCell "D37" "6.11"
Set-TextCell "E37" "  -0.53%  "

Set-TextCell "E39" "  +2.77%  "

Set-TextCell "D40" "1.00"
Set-TextCell "E40" "  -0.01%  "

Set-TextCell "D41" "176.92"
Set-TextCell "E41" "  +1.46%  "

Set-TextCell "E42" "  +1.60%  "

Set-TextCell "E43" "  -0.34%  "

Set-TextCell "D44" "47.04"
Set-TextCell "E44" "  -1.02%  "

Set-TextCell "D45" "2.76"
Set-TextCell "E45" "  +2.72%  "

Set-TextCell "E46" "  +2.09%  "

Set-TextCell "D47" "27.66"
Set-TextCell "E47" "  +0.04%  "

Set-TextCell "D48" "0.000273"
Set-TextCell "E48" "  -0.35%  "

Set-TextCell "E49" "  +0.60%  "

Set-TextCell "D50" "7.89"
Set-TextCell "E50" "  +1.70%  "

Set-TextCell "D51" "366.38"
Set-TextCell "E51" "  +2.11%  "
